$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 9).Value = '2025-05-19 17:35:18'
$ws.Cells.Item(28, 10).Value = '149.132.26.73'
$ws.Cells.Item(28, 11).Value = 'DT'
$ws.Cells.Item(28, 15).Value = 0
$ws.Cells.Item(28, 18).Value = 0
$ws.Cells.Item(28, 19).Value = 0
$ws.Cells.Item(28, 20).Value = 0
$ws.Cells.Item(28, 21).Value = 0
$ws.Cells.Item(28, 22).Value = 0
$ws.Cells.Item(28, 23).Value = 0
$ws.Cells.Item(28, 24).Value = 0
$ws.Cells.Item(28, 25).Value = '14% Malignant'
# Row 29
$ws.Cells.Item(29, 5).Value = 0
$ws.Cells.Item(29, 9).Value = '2025-05-19 17:50:34'
$ws.Cells.Item(29, 10).Value = '149.132.26.73'
$ws.Cells.Item(29, 11).Value = 'DT'
$ws.Cells.Item(29, 15).Value = 0
$ws.Cells.Item(29, 18).Value = 0
$ws.Cells.Item(29, 19).Value = 0
$ws.Cells.Item(29, 20).Value = 0
$ws.Cells.Item(29, 21).Value = 0
$ws.Cells.Item(29, 22).Value = 0
$ws.Cells.Item(29, 23).Value = 0
$ws.Cells.Item(29, 24).Value = 0
$ws.Cells.Item(29, 25).Value = '24% Malignant'
# Row 30
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(30, 9).Value = '2025-05-19 17:53:19'
$ws.Cells.Item(30, 10).Value = '149.132.26.73'
$ws.Cells.Item(30, 11).Value = 'DT'
$ws.Cells.Item(30, 15).Value = 0
$ws.Cells.Item(30, 18).Value = 0
$ws.Cells.Item(30, 19).Value = 0
$ws.Cells.Item(30, 20).Value = 0
$ws.Cells.Item(30, 21).Value = 0
$ws.Cells.Item(30, 22).Value = 0
$ws.Cells.Item(30, 23).Value = 0
$ws.Cells.Item(30, 24).Value = 0
$ws.Cells.Item(30, 25).Value = '24% Malignant'
# Row 31
$ws.Cells.Item(31, 5).Value = 1
$ws.Cells.Item(31, 9).Value = '2025-05-19 17:53:30'
$ws.Cells.Item(31, 10).Value = '149.132.26.73'
$ws.Cells.Item(31, 11).Value = 'DT'
$ws.Cells.Item(31, 15).Value = 0
$ws.Cells.Item(31, 18).Value = 0
$ws.Cells.Item(31, 19).Value = 0
$ws.Cells.Item(31, 20).Value = 1
$ws.Cells.Item(31, 21).Value = 0
$ws.Cells.Item(31, 22).Value = 1
$ws.Cells.Item(31, 23).Value = 0
$ws.Cells.Item(31, 24).Value = 0
$ws.Cells.Item(31, 25).Value = '100% Malignant'
# Row 32
$ws.Cells.Item(32, 5).Value = 1
$ws.Cells.Item(32, 9).Value = '2025-05-19 17:53:47'
$ws.Cells.Item(32, 10).Value = '149.132.26.73'
$ws.Cells.Item(32, 11).Value = 'DT'
$ws.Cells.Item(32, 15).Value = 0
$ws.Cells.Item(32, 18).Value = 0
$ws.Cells.Item(32, 19).Value = 0
$ws.Cells.Item(32, 20).Value = 1
$ws.Cells.Item(32, 21).Value = 0
$ws.Cells.Item(32, 22).Value = 1
$ws.Cells.Item(32, 23).Value = 1
$ws.Cells.Item(32, 24).Value = 0
$ws.Cells.Item(32, 25).Value = '100% Malignant'
# Row 33
$ws.Cells.Item(33, 5).Value = 0
$ws.Cells.Item(33, 9).Value = '2025-05-19 17:57:28'
$ws.Cells.Item(33, 10).Value = '149.132.26.73'
$ws.Cells.Item(33, 11).Value = 'DT'
$ws.Cells.Item(33, 15).Value = 0
$ws.Cells.Item(33, 18).Value = 0
$ws.Cells.Item(33, 19).Value = 0
$ws.Cells.Item(33, 20).Value = 0
$ws.Cells.Item(33, 21).Value = 0
$ws.Cells.Item(33, 22).Value = 0
$ws.Cells.Item(33, 23).Value = 0
$ws.Cells.Item(33, 24).Value = 0
$ws.Cells.Item(33, 25).Value = '24% Malignant'
# Row 34
$ws.Cells.Item(34, 5).Value = 1
$ws.Cells.Item(34, 9).Value = '2025-05-19 17:57:40'
$ws.Cells.Item(34, 10).Value = '149.132.26.73'
$ws.Cells.Item(34, 11).Value = 'DT'
$ws.Cells.Item(34, 15).Value = 0
$ws.Cells.Item(34, 18).Value = 0
$ws.Cells.Item(34, 19).Value = 0
$ws.Cells.Item(34, 20).Value = 1
$ws.Cells.Item(34, 21).Value = 1
$ws.Cells.Item(34, 22).Value = 1
$ws.Cells.Item(34, 23).Value = 0
$ws.Cells.Item(34, 24).Value = 0
$ws.Cells.Item(34, 25).Value = '67% Malignant'
# Row 35
$ws.Cells.Item(35, 5).Value = 1
$ws.Cells.Item(35, 9).Value = '2025-05-19 17:57:52'
$ws.Cells.Item(35, 10).Value = '149.132.26.73'
$ws.Cells.Item(35, 11).Value = 'DT'
$ws.Cells.Item(35, 15).Value = 0
$ws.Cells.Item(35, 18).Value = 0
$ws.Cells.Item(35, 19).Value = 0
$ws.Cells.Item(35, 20).Value = 1
$ws.Cells.Item(35, 21).Value = 0
$ws.Cells.Item(35, 22).Value = 1
$ws.Cells.Item(35, 23).Value = 0
$ws.Cells.Item(35, 24).Value = 0
$ws.Cells.Item(35, 25).Value = '100% Malignant'
# Row 36
$ws.Cells.Item(36, 5).Value = 0
$ws.Cells.Item(36, 9).Value = '2025-05-19 17:58:18'
$ws.Cells.Item(36, 10).Value = '149.132.26.73'
$ws.Cells.Item(36, 11).Value = 'DT'
$ws.Cells.Item(36, 15).Value = 0
$ws.Cells.Item(36, 18).Value = 0
$ws.Cells.Item(36, 19).Value = 0
$ws.Cells.Item(36, 20).Value = 0
$ws.Cells.Item(36, 21).Value = 0
$ws.Cells.Item(36, 22).Value = 1
$ws.Cells.Item(36, 23).Value = 1
$ws.Cells.Item(36, 24).Value = 0
$ws.Cells.Item(36, 25).Value = '29% Malignant'
# Row 37
$ws.Cells.Item(37, 5).Value = 0
$ws.Cells.Item(37, 9).Value = '2025-05-19 17:58:28'
$ws.Cells.Item(37, 10).Value = '149.132.26.73'
$ws.Cells.Item(37, 11).Value = 'DT'
$ws.Cells.Item(37, 15).Value = 0
$ws.Cells.Item(37, 18).Value = 0
$ws.Cells.Item(37, 19).Value = 0
$ws.Cells.Item(37, 20).Value = 0
$ws.Cells.Item(37, 21).Value = 0
$ws.Cells.Item(37, 22).Value = 0
$ws.Cells.Item(37, 23).Value = 1
$ws.Cells.Item(37, 24).Value = 0
$ws.Cells.Item(37, 25).Value = '0% Malignant'
# Row 38
$ws.Cells.Item(38, 5).Value = 1
$ws.Cells.Item(38, 9).Value = '2025-05-19 17:58:37'
$ws.Cells.Item(38, 10).Value = '149.132.26.73'
$ws.Cells.Item(38, 11).Value = 'DT'
$ws.Cells.Item(38, 15).Value = 0
$ws.Cells.Item(38, 18).Value = 0
$ws.Cells.Item(38, 19).Value = 0
$ws.Cells.Item(38, 20).Value = 1
$ws.Cells.Item(38, 21).Value = 1
$ws.Cells.Item(38, 22).Value = 0
$ws.Cells.Item(38, 23).Value = 0
$ws.Cells.Item(38, 24).Value = 0
$ws.Cells.Item(38, 25).Value = '84% Malignant'
# Row 39
$ws.Cells.Item(39, 5).Value = 1
$ws.Cells.Item(39, 9).Value = '2025-05-19 17:58:55'
$ws.Cells.Item(39, 10).Value = '149.132.26.73'
$ws.Cells.Item(39, 11).Value = 'DT'
$ws.Cells.Item(39, 15).Value = 0
$ws.Cells.Item(39, 18).Value = 0
$ws.Cells.Item(39, 19).Value = 0
$ws.Cells.Item(39, 20).Value = 1
$ws.Cells.Item(39, 21).Value = 1
$ws.Cells.Item(39, 22).Value = 1
$ws.Cells.Item(39, 23).Value = 0
$ws.Cells.Item(39, 24).Value = 0
$ws.Cells.Item(39, 25).Value = '67% Malignant'
# Row 40
$ws.Cells.Item(40, 5).Value = 0
$ws.Cells.Item(40, 9).Value = '2025-05-19 17:59:06'
$ws.Cells.Item(40, 10).Value = '149.132.26.73'
$ws.Cells.Item(40, 11).Value = 'DT'
$ws.Cells.Item(40, 15).Value = 0
$ws.Cells.Item(40, 18).Value = 0
$ws.Cells.Item(40, 19).Value = 0
$ws.Cells.Item(40, 20).Value = 1
$ws.Cells.Item(40, 21).Value = 1
$ws.Cells.Item(40, 22).Value = 1
$ws.Cells.Item(40, 23).Value = 0
$ws.Cells.Item(40, 24).Value = 0
$ws.Cells.Item(40, 25).Value = '50% Malignant'
# Row 41
$ws.Cells.Item(41, 5).Value = 1
$ws.Cells.Item(41, 9).Value = '2025-05-19 17:59:17'
$ws.Cells.Item(41, 10).Value = '149.132.26.73'
$ws.Cells.Item(41, 11).Value = 'DT'
$ws.Cells.Item(41, 15).Value = 0
$ws.Cells.Item(41, 18).Value = 0
$ws.Cells.Item(41, 19).Value = 0
$ws.Cells.Item(41, 20).Value = 1
$ws.Cells.Item(41, 21).Value = 1
$ws.Cells.Item(41, 22).Value = 1
$ws.Cells.Item(41, 23).Value = 0
$ws.Cells.Item(41, 24).Value = 0
$ws.Cells.Item(41, 25).Value = '67% Malignant'
# Row 42
$ws.Cells.Item(42, 5).Value = 1
$ws.Cells.Item(42, 9).Value = '2025-05-19 17:59:27'
$ws.Cells.Item(42, 10).Value = '149.132.26.73'
$ws.Cells.Item(42, 11).Value = 'DT'
$ws.Cells.Item(42, 15).Value = 0
$ws.Cells.Item(42, 18).Value = 0
$ws.Cells.Item(42, 19).Value = 0
$ws.Cells.Item(42, 20).Value = 1
$ws.Cells.Item(42, 21).Value = 1
$ws.Cells.Item(42, 22).Value = 1
$ws.Cells.Item(42, 23).Value = 1
$ws.Cells.Item(42, 24).Value = 0
$ws.Cells.Item(42, 25).Value = '100% Malignant'
# Row 43
$ws.Cells.Item(43, 5).Value = 1
$ws.Cells.Item(43, 9).Value = '2025-05-19 17:59:39'
$ws.Cells.Item(43, 10).Value = '149.132.26.73'
$ws.Cells.Item(43, 11).Value = 'DT'
$ws.Cells.Item(43, 15).Value = 0
$ws.Cells.Item(43, 18).Value = 0
$ws.Cells.Item(43, 19).Value = 0
$ws.Cells.Item(43, 20).Value = 1
$ws.Cells.Item(43, 21).Value = 1
$ws.Cells.Item(43, 22).Value = 0
$ws.Cells.Item(43, 23).Value = 1
$ws.Cells.Item(43, 24).Value = 0
$ws.Cells.Item(43, 25).Value = '75% Malignant'
